# Auto-generated: updates currentAveragePrice / LevePrice / LeveProfit columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 374.72726
$ws.Range("J2").Value = 790.75
$ws.Range("L2").Value = 790.75
$ws.Range("N2").Value = -1016.75
$ws.Range("H40").Value = 11118901
$ws.Range("I40").Value = 4474.5
$ws.Range("J40").Value = 18528518
$ws.Range("K40").Value = 4474.5
$ws.Range("L40").Value = 18528518
$ws.Range("M40").Value = -4299.5
$ws.Range("N40").Value = -18528868
$ws.Range("H43").Value = 6131.8335
$ws.Range("I43").Value = 7000
$ws.Range("J43").Value = 5697.75
$ws.Range("K43").Value = 7000
$ws.Range("L43").Value = 5697.75
$ws.Range("M43").Value = -6931
$ws.Range("N43").Value = -5835.75
$ws.Range("H64").Value = 25003922
$ws.Range("I64").Value = 28575196
$ws.Range("K64").Value = 28575196
$ws.Range("M64").Value = -28574948
$ws.Range("H67").Value = 25003922
$ws.Range("I67").Value = 28575196
$ws.Range("K67").Value = 28575196
$ws.Range("M67").Value = -28574338
$ws.Range("H69").Value = 16778.934
$ws.Range("I69").Value = 7844.5
$ws.Range("J69").Value = 18153.46
$ws.Range("K69").Value = 23533.5
$ws.Range("L69").Value = 54460.38
$ws.Range("M69").Value = -22659.5
$ws.Range("N69").Value = -56208.38
$ws.Range("H72").Value = 16778.934
$ws.Range("I72").Value = 7844.5
$ws.Range("J72").Value = 18153.46
$ws.Range("K72").Value = 70600.5
$ws.Range("L72").Value = 163381.14
$ws.Range("M72").Value = -66232.5
$ws.Range("N72").Value = -172117.14
$ws.Range("H74").Value = 55563556
$ws.Range("I74").Value = 55563556
$ws.Range("K74").Value = 55563556
$ws.Range("M74").Value = -55562620
$ws.Range("H76").Value = 3824.25
$ws.Range("I76").Value = 3418.8
$ws.Range("K76").Value = 3418.8
$ws.Range("M76").Value = -3103.8
$ws.Range("H77").Value = 55563556
$ws.Range("I77").Value = 55563556
$ws.Range("K77").Value = 277817780
$ws.Range("M77").Value = -277813100
$ws.Range("H79").Value = 3824.25
$ws.Range("I79").Value = 3418.8
$ws.Range("K79").Value = 3418.8
$ws.Range("M79").Value = -2326.8
$ws.Range("H112").Value = 102858.85
$ws.Range("I112").Value = 251499.5
$ws.Range("J112").Value = 65698.69
$ws.Range("K112").Value = 754498.5
$ws.Range("L112").Value = 197096.07
$ws.Range("M112").Value = -753390.5
$ws.Range("N112").Value = -199312.07
$ws.Range("H132").Value = 2406.842
$ws.Range("I132").Value = 2639.375
$ws.Range("J132").Value = 1166.6666
$ws.Range("K132").Value = 7918.125
$ws.Range("L132").Value = 3499.9998
$ws.Range("M132").Value = -5388.125
$ws.Range("N132").Value = -8559.9998
$ws.Range("H137").Value = 3860.484
$ws.Range("I137").Value = 4406.6
$ws.Range("J137").Value = 3348.5
$ws.Range("K137").Value = 13219.8
$ws.Range("L137").Value = 10045.5
$ws.Range("M137").Value = -10669.8
$ws.Range("N137").Value = -15145.5
$ws.Range("H138").Value = 6979.22
$ws.Range("J138").Value = 8362.275
$ws.Range("L138").Value = 25086.825
$ws.Range("N138").Value = -35366.825

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 769.7059
$ws.Range("I5").Value = 705.6667
$ws.Range("J5").Value = 1250
$ws.Range("K5").Value = 705.6667
$ws.Range("L5").Value = 1250
$ws.Range("M5").Value = -593.6667
$ws.Range("N5").Value = -1474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 769.7059
$ws.Range("I4").Value = 705.6667
$ws.Range("J4").Value = 1250
$ws.Range("K4").Value = 705.6667
$ws.Range("L4").Value = 1250
$ws.Range("M4").Value = -590.6667
$ws.Range("N4").Value = -1480
$ws.Range("H86").Value = 3303.4614
$ws.Range("I86").Value = 3322.3635
$ws.Range("J86").Value = 3199.5
$ws.Range("K86").Value = 3322.3635
$ws.Range("L86").Value = 3199.5
$ws.Range("M86").Value = -2199.3635
$ws.Range("N86").Value = -5445.5
$ws.Range("H89").Value = 3303.4614
$ws.Range("I89").Value = 3322.3635
$ws.Range("J89").Value = 3199.5
$ws.Range("K89").Value = 16611.8175
$ws.Range("L89").Value = 15997.5
$ws.Range("M89").Value = -10995.8175
$ws.Range("N89").Value = -27229.5
$ws.Range("H105").Value = 3226.3333
$ws.Range("I105").Value = 2626
$ws.Range("J105").Value = 4126.8335
$ws.Range("K105").Value = 2626
$ws.Range("L105").Value = 4126.8335
$ws.Range("M105").Value = -879
$ws.Range("N105").Value = -7620.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 212.71428
$ws.Range("I7").Value = 322.75
$ws.Range("J7").Value = 66
$ws.Range("K7").Value = 322.75
$ws.Range("L7").Value = 66
$ws.Range("M7").Value = -209.75
$ws.Range("N7").Value = -292
$ws.Range("H22").Value = 3755.3103
$ws.Range("I22").Value = 3755.3103
$ws.Range("K22").Value = 3755.3103
$ws.Range("M22").Value = -3405.3103
$ws.Range("H31").Value = 5181.8623
$ws.Range("I31").Value = 3188.5417
$ws.Range("K31").Value = 3188.5417
$ws.Range("M31").Value = -2893.5417
$ws.Range("H34").Value = 5181.8623
$ws.Range("I34").Value = 3188.5417
$ws.Range("K34").Value = 3188.5417
$ws.Range("M34").Value = -2986.5417
$ws.Range("H62").Value = 1966.3334
$ws.Range("I62").Value = 1966.3334
$ws.Range("K62").Value = 1966.3334
$ws.Range("M62").Value = -1342.3334
$ws.Range("H65").Value = 1966.3334
$ws.Range("I65").Value = 1966.3334
$ws.Range("K65").Value = 9831.666999999999
$ws.Range("M65").Value = -6711.666999999999
$ws.Range("H97").Value = 38999
$ws.Range("J97").Value = 38999
$ws.Range("L97").Value = 38999
$ws.Range("N97").Value = -40981
$ws.Range("H122").Value = 3364
$ws.Range("I122").Value = 3256.6428
$ws.Range("J122").Value = 3614.5
$ws.Range("K122").Value = 9769.928400000001
$ws.Range("L122").Value = 10843.5
$ws.Range("M122").Value = -7319.928400000001
$ws.Range("N122").Value = -15743.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1040.6
$ws.Range("I2").Value = 32.6
$ws.Range("J2").Value = 2048.6
$ws.Range("K2").Value = 195.6
$ws.Range("L2").Value = 12291.6
$ws.Range("M2").Value = -82.60000000000002
$ws.Range("N2").Value = -12517.6
$ws.Range("H23").Value = 574.5
$ws.Range("I23").Value = 458.57144
$ws.Range("J23").Value = 690.4286
$ws.Range("K23").Value = 1375.71432
$ws.Range("L23").Value = 2071.2858
$ws.Range("M23").Value = -1140.71432
$ws.Range("N23").Value = -2541.2858
$ws.Range("H38").Value = 178.2
$ws.Range("J38").Value = 97.333336
$ws.Range("L38").Value = 292.000008
$ws.Range("N38").Value = -986.000008
$ws.Range("H59").Value = 5549
$ws.Range("I59").Value = 99
$ws.Range("K59").Value = 297
$ws.Range("M59").Value = 243
$ws.Range("H107").Value = 981.4
$ws.Range("I107").Value = 205.2
$ws.Range("J107").Value = 1369.5
$ws.Range("K107").Value = 615.5999999999999
$ws.Range("L107").Value = 4108.5
$ws.Range("M107").Value = 1304.4
$ws.Range("N107").Value = -7948.5
$ws.Range("H113").Value = 167027.67
$ws.Range("J113").Value = 444.25
$ws.Range("L113").Value = 1332.75
$ws.Range("N113").Value = -5672.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 50499.375
$ws.Range("J15").Value = 50499.375
$ws.Range("L15").Value = 50499.375
$ws.Range("N15").Value = -51075.375
$ws.Range("H54").Value = 19214.666
$ws.Range("J54").Value = 19214.666
$ws.Range("L54").Value = 19214.666
$ws.Range("N54").Value = -19994.666
$ws.Range("H81").Value = 50499.375
$ws.Range("J81").Value = 50499.375
$ws.Range("L81").Value = 50499.375
$ws.Range("N81").Value = -52495.375
$ws.Range("H84").Value = 50499.375
$ws.Range("J84").Value = 50499.375
$ws.Range("L84").Value = 151498.125
$ws.Range("N84").Value = -161482.125
$ws.Range("H92").Value = 7419.25
$ws.Range("J92").Value = 7419.25
$ws.Range("L92").Value = 7419.25
$ws.Range("N92").Value = -11163.25
$ws.Range("H132").Value = 2671810.2
$ws.Range("I132").Value = 2852783
$ws.Range("J132").Value = 17541.666
$ws.Range("K132").Value = 8558349
$ws.Range("L132").Value = 52624.99800000001
$ws.Range("M132").Value = -8555819
$ws.Range("N132").Value = -57684.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3187.5
$ws.Range("I22").Value = 3195.5715
$ws.Range("K22").Value = 3195.5715
$ws.Range("M22").Value = -2900.5715
$ws.Range("H27").Value = 3187.5
$ws.Range("I27").Value = 3195.5715
$ws.Range("K27").Value = 3195.5715
$ws.Range("M27").Value = -3088.5715
$ws.Range("H40").Value = 3766.3333
$ws.Range("I40").Value = 3766.3333
$ws.Range("K40").Value = 3766.3333
$ws.Range("M40").Value = -3630.3333
$ws.Range("H46").Value = 1682.5714
$ws.Range("I46").Value = 1739.2
$ws.Range("K46").Value = 1739.2
$ws.Range("M46").Value = -1551.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 29414826
$ws.Range("I132").Value = 50002836
$ws.Range("K132").Value = 150008508
$ws.Range("M132").Value = -150005978
$ws.Range("H136").Value = 20001896
$ws.Range("I136").Value = 20835162
$ws.Range("K136").Value = 62505486
$ws.Range("M136").Value = -62502936
